# Update "想去人数" (F column) figures (and one "最低票价" / G13 sold-out flag)
# to match the newly scraped data (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 1356
$ws1.Range("F7").Value  = 623
$ws1.Range("F9").Value  = 74
$ws1.Range("F10").Value = 430
$ws1.Range("F13").Value = 30630
$ws1.Range("G13").Value = "已售罄"
$ws1.Range("F14").Value = 6211
$ws1.Range("F16").Value = 311
$ws1.Range("F18").Value = 86
$ws1.Range("F20").Value = 63
$ws1.Range("F22").Value = 395
$ws1.Range("F24").Value = 722
$ws1.Range("F30").Value = 27
$ws1.Range("F31").Value = 701
$ws1.Range("F32").Value = 250
$ws1.Range("F33").Value = 118
$ws1.Range("F34").Value = 644
$ws1.Range("F35").Value = 90
$ws1.Range("F37").Value = 719
$ws1.Range("F40").Value = 21

# --- Sheet 2: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value  = 1038
$ws2.Range("F5").Value  = 118
$ws2.Range("F6").Value  = 283
$ws2.Range("F15").Value = 38

# --- Sheet 3: 本地生活 ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F5").Value = 323

# --- Sheet 4: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 323
$ws4.Range("F6").Value  = 1038
$ws4.Range("F9").Value  = 1356
$ws4.Range("F11").Value = 623
$ws4.Range("F12").Value = 74
$ws4.Range("F13").Value = 430
$ws4.Range("F17").Value = 118
$ws4.Range("F18").Value = 118
$ws4.Range("F19").Value = 283
$ws4.Range("F24").Value = 311
$ws4.Range("F27").Value = 86
$ws4.Range("F29").Value = 63
$ws4.Range("F32").Value = 38
$ws4.Range("F33").Value = 395
$ws4.Range("F35").Value = 722
$ws4.Range("F41").Value = 27
$ws4.Range("F42").Value = 701
$ws4.Range("F44").Value = 250
$ws4.Range("F45").Value = 118
$ws4.Range("F46").Value = 90
$ws4.Range("F47").Value = 719
